$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 900
$ws.Range("I4").Value = 900
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -786

$ws.Range("H6").Value = 2600.5
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 2600.5
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 7801.5
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -8025.5

$ws.Range("H12").Value = 1033
$ws.Range("I12").Value = 880
$ws.Range("J12").Value = 1186
$ws.Range("K12").Value = 880
$ws.Range("L12").Value = 1186
$ws.Range("M12").Value = -710
$ws.Range("N12").Value = -1526

$ws.Range("H74").Value = 4198.2856
$ws.Range("I74").Value = 4198
$ws.Range("J74").Value = 4198.5
$ws.Range("K74").Value = 4198
$ws.Range("L74").Value = 4198.5
$ws.Range("M74").Value = -3262
$ws.Range("N74").Value = -6070.5

$ws.Range("H77").Value = 4198.2856
$ws.Range("I77").Value = 4198
$ws.Range("J77").Value = 4198.5
$ws.Range("K77").Value = 20990
$ws.Range("L77").Value = 20992.5
$ws.Range("M77").Value = -16310
$ws.Range("N77").Value = -30352.5

$ws.Range("H88").Value = 1814.7273
$ws.Range("I88").Value = 8999
$ws.Range("J88").Value = 1096.3
$ws.Range("K88").Value = 8999
$ws.Range("L88").Value = 1096.3
$ws.Range("M88").Value = -8593
$ws.Range("N88").Value = -1908.3

$ws.Range("H91").Value = 1814.7273
$ws.Range("I91").Value = 8999
$ws.Range("J91").Value = 1096.3
$ws.Range("K91").Value = 8999
$ws.Range("L91").Value = 1096.3
$ws.Range("M91").Value = -7595
$ws.Range("N91").Value = -3904.3

$ws.Range("H112").Value = 3799.1428
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 4060.6155
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 12181.8465
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -14397.8465

$ws.Range("H132").Value = 1002.87177
$ws.Range("I132").Value = 1008.2105
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 3024.6315
$ws.Range("L132").Value = 2400
$ws.Range("M132").Value = -494.6315
$ws.Range("N132").Value = -7460

$ws.Range("H138").Value = 2675.1562
$ws.Range("I138").Value = 3495
$ws.Range("J138").Value = 1746
$ws.Range("K138").Value = 10485
$ws.Range("L138").Value = 5238
$ws.Range("M138").Value = -5345
$ws.Range("N138").Value = -15518

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3293.6667
$ws.Range("I32").Value = 2120.7737
$ws.Range("J32").Value = 8075.4614
$ws.Range("K32").Value = 2120.7737
$ws.Range("L32").Value = 8075.4614
$ws.Range("M32").Value = -1833.7737
$ws.Range("N32").Value = -8649.4614

$ws.Range("H33").Value = 12222
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 12222
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 12222
$ws.Range("N33").Value = -12880

$ws.Range("H61").Value = 2985.375
$ws.Range("I61").Value = 1570.8462
$ws.Range("J61").Value = 9115
$ws.Range("K61").Value = 1570.8462
$ws.Range("L61").Value = 9115
$ws.Range("M61").Value = -1358.8462
$ws.Range("N61").Value = -9539

$ws.Range("H112").Value = 40000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 40000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

$ws.Range("H122").Value = 1610.4667
$ws.Range("I122").Value = 1598.2858
$ws.Range("J122").Value = 1621.125
$ws.Range("K122").Value = 4794.857400000001
$ws.Range("L122").Value = 4863.375
$ws.Range("M122").Value = -2344.857400000001
$ws.Range("N122").Value = -9763.375

$ws.Range("H132").Value = 2523.3333
$ws.Range("I132").Value = 1816
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 5448
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2918
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 2985.375
$ws.Range("I136").Value = 1570.8462
$ws.Range("J136").Value = 9115
$ws.Range("K136").Value = 4712.5386
$ws.Range("L136").Value = 27345
$ws.Range("M136").Value = -2162.5386
$ws.Range("N136").Value = -32445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8159.4
$ws.Range("I134").Value = 10009.521
$ws.Range("J134").Value = 4613.3335
$ws.Range("K134").Value = 30028.563
$ws.Range("L134").Value = 13840.0005
$ws.Range("M134").Value = -27493.563
$ws.Range("N134").Value = -18910.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 862.7143
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 939.8333
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 939.8333
$ws.Range("M22").Value = -50
$ws.Range("N22").Value = -1639.8333

$ws.Range("H68").Value = 45000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 45000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 45000
$ws.Range("N68").Value = -46498

$ws.Range("H71").Value = 45000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 45000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 135000
$ws.Range("N71").Value = -142488

$ws.Range("H129").Value = 44999.25
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 44999.25
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 44999.25
$ws.Range("N129").Value = -54999.25

$ws.Range("H132").Value = 2108.4614
$ws.Range("I132").Value = 1124
$ws.Range("J132").Value = 5390
$ws.Range("K132").Value = 3372
$ws.Range("L132").Value = 16170
$ws.Range("M132").Value = -842
$ws.Range("N132").Value = -21230

$ws.Range("H134").Value = 784.35297
$ws.Range("I134").Value = 708.375
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2125.125
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 409.875
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 799
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 799
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 2397
$ws.Range("N98").Value = -5393

$ws.Range("H114").Value = 2171.8
$ws.Range("I114").Value = 1014
$ws.Range("J114").Value = 2943.6667
$ws.Range("K114").Value = 3042
$ws.Range("L114").Value = 8831.000100000001
$ws.Range("M114").Value = 212
$ws.Range("N114").Value = -15339.0001

$ws.Range("H118").Value = 2284.5
$ws.Range("I118").Value = 569
$ws.Range("J118").Value = 4000
$ws.Range("K118").Value = 1707
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = -464
$ws.Range("N118").Value = -14486

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("M128").ClearContents()

$ws.Range("H132").Value = 1600
$ws.Range("I132").Value = 1600
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14400
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11870
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 41000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 41000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 41000
$ws.Range("N42").Value = -41970

$ws.Range("H102").Value = 2058.4517
$ws.Range("I102").Value = 1980.9
$ws.Range("J102").Value = 2199.4546
$ws.Range("K102").Value = 1980.9
$ws.Range("L102").Value = 2199.4546
$ws.Range("M102").Value = -358.9000000000001
$ws.Range("N102").Value = -5443.4546

$ws.Range("H115").Value = 41000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 41000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 41000
$ws.Range("N115").Value = -43350

$ws.Range("H125").Value = 30000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 30000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920

$ws.Range("H132").Value = 5255.75
$ws.Range("I132").Value = 4252.647
$ws.Range("J132").Value = 7691.857
$ws.Range("K132").Value = 12757.941
$ws.Range("L132").Value = 23075.571
$ws.Range("M132").Value = -10227.941
$ws.Range("N132").Value = -28135.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 654.6667
$ws.Range("I93").Value = 232.5
$ws.Range("J93").Value = 1499
$ws.Range("K93").Value = 232.5
$ws.Range("L93").Value = 1499
$ws.Range("M93").Value = 1015.5
$ws.Range("N93").Value = -3995

$ws.Range("H110").Value = 13547.667
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 13547.667
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 13547.667
$ws.Range("N110").Value = -21727.667

$ws.Range("H122").Value = 10751.25
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 13335
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 40005
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -44905

$ws.Range("H132").Value = 1979
$ws.Range("I132").Value = 1125
$ws.Range("J132").Value = 2168.7778
$ws.Range("K132").Value = 3375
$ws.Range("L132").Value = 6506.3334
$ws.Range("M132").Value = -845
$ws.Range("N132").Value = -11566.3334

$ws.Range("H136").Value = 3396.7273
$ws.Range("I136").Value = 1929.4
$ws.Range("J136").Value = 4619.5
$ws.Range("K136").Value = 5788.200000000001
$ws.Range("L136").Value = 13858.5
$ws.Range("M136").Value = -3238.200000000001
$ws.Range("N136").Value = -18958.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 46779.2
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 46779.2
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 46779.2
$ws.Range("N123").Value = -56579.2

$ws.Range("H132").Value = 3688.9443
$ws.Range("I132").Value = 3426.9333
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 10280.7999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -7750.7999
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 3967.0833
$ws.Range("I136").Value = 3901.3333
$ws.Range("J136").Value = 3989
$ws.Range("K136").Value = 11703.9999
$ws.Range("L136").Value = 11967
$ws.Range("M136").Value = -9153.999899999999
$ws.Range("N136").Value = -17067
